$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.575.73"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").Value = "2.043.53"
$ws.Range("E3").Value = "  +3.35%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.94"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("E11").Value = "  +1.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.93%  "

$ws.Range("D14").Value = "2.344.11"
$ws.Range("E14").Value = "  +3.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.827"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.04%  "

$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").Value = "2.045.61"
$ws.Range("E18").Value = "  +3.58%  "

$ws.Range("D19").Value = "37.416.53"
$ws.Range("E19").Value = "  +2.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.50%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0670"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.24%  "

$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("E44").Value = "  +3.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.81%  "

$ws.Range("D46").Value = "1.406.83"
$ws.Range("E46").Value = "  +2.75%  "

$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("E48").Value = "  +2.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("E50").Value = "  +8.40%  "

$ws.Range("E51").Value = "  +2.36%  "
